$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 3, 4 and 5 get their non-date columns cyclically shifted:
#   new row 3 <- old row 4 (Calidad..Precio $/Kg)
#   new row 4 <- old row 5 (Calidad..Precio $/Kg)
#   new row 5 <- old row 3 (Calidad..Precio $/Kg)
# Dates (column D) swap between row 3 and row 5; row 4's date is unchanged.

# Row 3 (date changes: 44316 -> 44280)
$ws.Range("D3").Value = 44280
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("S3").Value = 806

# Row 4 (date unchanged = 44280)
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("S4").Value = 667

# Row 5 (date changes: 44280 -> 44316)
$ws.Range("D5").Value = 44316
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("S5").Value = 1111
